$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated "Price" cells (column D) contain values that look like plain
# numbers (e.g. "218.28", "10.10"). Column D is free-form text in this sheet
# (prices are sometimes formatted like "26.149.39" which is not a valid
# number), so force those specific cells to Text format first to keep Excel
# from silently re-interpreting the new value as a number and dropping
# significant trailing/leading zeros.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D12", "D14", "D19", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D45", "D47", "D48", "D49", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.149.39'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '1.657.44'
$ws.Range("E3").Value = '  -0.35%  '
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("D5").Value = '218.28'
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").Value = '0.5289'
$ws.Range("E6").Value = '  +1.10%  '
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").Value = '0.2613'
$ws.Range("E8").Value = '  -2.22%  '
$ws.Range("D9").Value = '0.06348'
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("D10").Value = '20.46'
$ws.Range("E10").Value = '  -1.67%  '
$ws.Range("D11").Value = '0.07790'
$ws.Range("E11").Value = '  +0.80%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '4.512'
$ws.Range("E12").Value = '  +1.76%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.656.19'
$ws.Range("E13").Value = '  -0.52%  '
$ws.Range("D14").Value = '0.5498'
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("D15").Value = '0.0₅8206'
$ws.Range("E15").Value = '  +0.35%  '
$ws.Range("E16").Value = '  +1.39%  '
$ws.Range("D17").Value = '26.148.85'
$ws.Range("E17").Value = '  -0.40%  '
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").Value = '4.582'
$ws.Range("E19").Value = '  -1.36%  '
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("D21").Value = '10.10'
$ws.Range("E21").Value = '  +0.35%  '
$ws.Range("D22").Value = '6.041'
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").Value = '141.86'
$ws.Range("E24").Value = '  +1.39%  '
$ws.Range("E25").Value = '  +1.68%  '
$ws.Range("D26").Value = '7.284'
$ws.Range("E26").Value = '  +1.62%  '
$ws.Range("D27").Value = '16.20'
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("D28").Value = '1.442'
$ws.Range("E28").Value = '  +2.01%  '
$ws.Range("E29").Value = '  -2.97%  '
$ws.Range("D30").Value = '1.281'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").Value = '3.526'
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("D32").Value = '3.268'
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("E33").Value = '  -2.11%  '
$ws.Range("D34").Value = '0.9558'
$ws.Range("E34").Value = '  -1.20%  '
$ws.Range("D35").Value = '2.794'
$ws.Range("D36").Value = '2.414'
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("D37").Value = '0.5704'
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("D38").Value = '0.01621'
$ws.Range("E38").Value = '  +1.76%  '
$ws.Range("D39").Value = '5.810'
$ws.Range("E39").Value = '  -2.90%  '
$ws.Range("D40").Value = '0.8488'
$ws.Range("E40").Value = '  -0.95%  '
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("D42").Value = '103.25'
$ws.Range("E42").Value = '  +2.96%  '
$ws.Range("D43").Value = '1.025.57'
$ws.Range("E43").Value = '  +1.05%  '
$ws.Range("D44").Value = '1.802.10'
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").Value = '57.42'
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("D47").Value = '1.487'
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("D48").Value = '0.4293'
$ws.Range("D49").Value = '0.05159'
$ws.Range("E49").Value = '  -0.49%  '
$ws.Range("E50").Value = '  -2.36%  '
$ws.Range("D51").Value = '0.09719'
$ws.Range("E51").Value = '  +0.02%  '
